$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(20).Insert()
$ws.Range("B20:C20").Borders.LineStyle = 1
